# Add a "2022-Q3" sheet (new fund-holdings data) right after "总计" and
# before "2022-Q2", and record its summary row on the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new first data row for
#    2022-Q3 and push the existing three rows down by one.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Row 5 is brand new (the sheet only had rows 1..4 before) - copy the
# index-column format from row 4 so A5 picks up the same style as the
# other index cells instead of defaulting to unstyled.
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)

$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(5,2).Value = "2021-Q4"
$total.Cells.Item(5,3).Value = 2
$total.Cells.Item(5,4).Value = 0.17

$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = "2022-Q1"
$total.Cells.Item(4,3).Value = 6
$total.Cells.Item(4,4).Value = 0.07000000000000001

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2022-Q2"
$total.Cells.Item(3,3).Value = 4
$total.Cells.Item(3,4).Value = 0.62

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 20
$total.Cells.Item(2,4).Value = 1.56

# ---------------------------------------------------------------------
# 2. Create the new "2022-Q3" sheet by copying the "2022-Q2" sheet's
#    layout/formatting, positioned right before it, then overwrite the
#    cell values with the 2022-Q3 fund-holdings data.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.ActiveSheet
$q3.Name = "2022-Q3"

# Fund rows for 2022-Q3 (idx, code, name, scale, totalPos, posPct, mv, rank)
$rows = @(
  @(0, "001743", "诺安优选回报灵活配置混合", "13.65", "73.02", "3.32", "0.4532", 7),
  @(1, "003501", "泰达宏利睿智稳健灵活配置混合A", "10.30", "76.72", "1.80", "0.1854", 10),
  @(2, "162204", "泰达宏利行业精选混合A", "8.57", "81.08", "1.76", "0.1508", 9),
  @(3, "011346", "淳厚鑫淳一年持有期混合", "4.81", "69.72", "2.69", "0.1294", 7),
  @(4, "013280", "泰达宏利睿智稳健灵活配置混合C", "5.95", "76.72", "1.80", "0.1071", 10),
  @(5, "162203", "泰达宏利稳定混合", "3.13", "91.16", "3.20", "0.1002", 9),
  @(6, "519013", "海富通风格优势混合", "3.41", "87.16", "2.87", "0.0979", 7),
  @(7, "015601", "泰达宏利行业精选混合C", "4.13", "81.08", "1.76", "0.0727", 9),
  @(8, "011431", "泰达宏利消费服务混合A", "1.47", "86.98", "3.98", "0.0585", 2),
  @(9, "012454", "淳厚鑫悦混合A", "2.06", "75.61", "2.62", "0.0540", 7),
  @(10, "005741", "南方君信灵活配置混合A", "2.80", "71.98", "1.92", "0.0538", 7),
  @(11, "519139", "海富通沪港深灵活配置混合", "0.71", "94.17", "3.64", "0.0258", 7),
  @(12, "012455", "淳厚鑫悦混合C", "0.68", "75.61", "2.62", "0.0178", 7),
  @(13, "011432", "泰达宏利消费服务混合C", "0.44", "86.98", "3.98", "0.0175", 2),
  @(14, "005493", "鑫元价值精选灵活配置混合A", "0.55", "76.82", "2.97", "0.0163", 8),
  @(15, "001744", "诺安进取回报灵活配置混合", "0.23", "82.31", "4.11", "0.0095", 6),
  @(16, "010150", "南方君信灵活配置混合C", "0.20", "71.98", "1.92", "0.0038", 7),
  @(17, "006193", "鑫元核心资产股票A", "0.11", "83.48", "3.46", "0.0038", 5),
  @(18, "005494", "鑫元价值精选灵活配置混合C", "0.01", "76.82", "2.97", "0.0003", 8),
  @(19, "006194", "鑫元核心资产股票C", "0.01", "83.48", "3.46", "0.0003", 5)
)

$lastExistingRow = 5  # rows 2..5 already exist (copied from 2022-Q2, 4 funds)
$lastNewRow = 1 + $rows.Count  # rows 2..21

# Grow the sheet: copy the formatting of row 2 onto any brand-new rows
# (rows 6..21) before writing values into them.
if ($lastNewRow -gt $lastExistingRow) {
    $q3.Range("A2:H2").Copy()
    $q3.Range("A" + ($lastExistingRow + 1) + ":H" + $lastNewRow).PasteSpecial(-4122)
}

foreach ($row in $rows) {
    $r = 2 + [int]$row[0]

    $q3.Cells.Item($r,1).Value = [int]$row[0]

    $q3.Cells.Item($r,2).Value = "'" + $row[1]
    $q3.Cells.Item($r,2).Style = "Normal"

    $q3.Cells.Item($r,3).Value = $row[2]

    $q3.Cells.Item($r,4).Value = "'" + $row[3]
    $q3.Cells.Item($r,4).Style = "Normal"

    $q3.Cells.Item($r,5).Value = "'" + $row[4]
    $q3.Cells.Item($r,5).Style = "Normal"

    $q3.Cells.Item($r,6).Value = "'" + $row[5]
    $q3.Cells.Item($r,6).Style = "Normal"

    $q3.Cells.Item($r,7).Value = "'" + $row[6]
    $q3.Cells.Item($r,7).Style = "Normal"

    $q3.Cells.Item($r,8).Value = [int]$row[7]
}
